# "seller nao vetex sem pagamento"
# Reposition four shapes vertically on slide 3 (the SKU / VTEXADMIN / SITE
# integration diagram): the "R$ $" price rectangle + its "$ $" label move
# down a bit, and the two "ATIVA PRODUTO" / "ATIVA SKU" arrows swap places
# (PRODUTO moves down past where SKU used to be, SKU moves up past where
# PRODUTO used to be).
#
# EMU -> points uses 914400 EMU/inch, 12700 EMU/point. PowerPoint's
# Shape.Top/.Left are single-precision floats in points, so a plain
# emu/12700.0 division can truncate one EMU low once it round-trips through
# that float32 storage. A tiny (2e-5 pt ~= 0.25 EMU) upward nudge keeps the
# intended value from flooring down without overshooting into the next EMU.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

function EmuToPt([double]$emu) {
    return ($emu / 12700.0) + 0.00002
}

# Rectangle 83 ("R$ $" fill rectangle): off x=1180685 y=3823405 -> y=3993224
$rect83 = $s.Shapes.Item("Rectangle 83")
$rect83.Top = EmuToPt 3993224

# TextBox 73 ("$ $" red label, rotated): off x=1130973 y=3761641 -> y=3931460
$textBox73 = $s.Shapes.Item("TextBox 73")
$textBox73.Top = EmuToPt 3931460

# Right Arrow (text "ATIVA PRODUTO"): off x=3020070 y=5175340 -> y=5575376
$ativaProduto = $s.Shapes.Item(49)
$ativaProduto.Top = EmuToPt 5575376

# Right Arrow (text "ATIVA SKU"): off x=3020070 y=5682338 -> y=4947578
$ativaSku = $s.Shapes.Item(50)
$ativaSku.Top = EmuToPt 4947578
